{"js": "// Update the address line in the resume header.\n// Old: \"1613 Black Duck Terrace Apt F, Carrollton, TX 75010\"\n// New: \"2006 Azure Pointe, Richardson TX 75080\"\nconst oldAddress = \"1613 Black Duck Terrace Apt F, Carrollton, TX 75010\";\nconst newAddress = \"2006 Azure Pointe, Richardson TX 75080\";\n\nconst body = context.document.body;\nconst results = body.search(oldAddress, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // Replace the matched range's text in place so run formatting (font,\n  // size, etc.) carried by the existing run is preserved.\n  results.items[i].insertText(newAddress, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the address line in the resume header.\n# Old: \"1613 Black Duck Terrace Apt F, Carrollton, TX 75010\"\n# New: \"2006 Azure Pointe, Richardson TX 75080\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"1613 Black Duck Terrace Apt F, Carrollton, TX 75010\"\n$find.Replacement.Text = \"2006 Azure Pointe, Richardson TX 75080\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\n  [ref]$find.Text,\n  [ref]$find.MatchCase,\n  [ref]$find.MatchWholeWord,\n  [ref]$find.MatchWildcards,\n  [ref]$false,\n  [ref]$false,\n  [ref]$find.Forward,\n  [ref]$find.Wrap,\n  [ref]$false,\n  [ref]$find.Replacement.Text,\n  [ref]2\n)\n"}
